$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5.4, 1, 1, 2, 3, 10, 10, 10, 10, 0.04, 0.03, 0.02, 0.07, 200, 300, 600, 0.282, 0.75),
    @(5.8, 1, 1, 2, 3, 10, 10, 10, 10, 0.04, 0.03, 0.02, 0.07, 200, 300, 600, 0.214, 0.76),
    @(5.8, 3, 1, 2, 3, 10, 10, 10, 10, 0.04, 0.03, 0.02, 0.07, 200, 300, 600, 0.072, 0.76)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

$ws.Range("A4:R4").Select()
